# Append, after the existing "aa" paragraph:
#   (empty paragraph)
#   "save it"
#   (empty paragraph)
#   "okay"
#
# Plain paragraph-insertion APIs (Range.InsertParagraphAfter,
# Paragraphs.Add, Selection.TypeParagraph, ...) always leave a stray
# empty run (<w:r/>) inside a "blank" paragraph. To get a truly empty
# <w:p/> (no run at all) - matching how Word collapses a paragraph that
# never had any content typed into it - insert raw WordprocessingML via
# Range.InsertXML at a zero-length Range positioned at the very end of
# the document's main story.

$d = $word.ActiveDocument

# A fresh Range object anchored at the end of the document content.
# (Using the Paragraph's own .Range + Collapse() instead would bind the
# insert to that paragraph's span and clobber its text, so we build the
# target Range straight from the document.)
$endPos = $d.Content.End
$insertionPoint = $d.Range($endPos, $endPos)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$newParagraphsXml = "<w:p $wNs/>" + `
    "<w:p $wNs><w:r><w:t>save it</w:t></w:r></w:p>" + `
    "<w:p $wNs/>" + `
    "<w:p $wNs><w:r><w:t>okay</w:t></w:r></w:p>"

$insertionPoint.InsertXML($newParagraphsXml)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
